{"js": "// Change \"Version 2.\" to \"Version 1.\" in the document body.\nconst body = context.document.body;\nconst results = body.search(\"Version 2.\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  results.items[0].insertText(\"Version 1.\", \"Replace\");\n} else {\n  // Fallback: search for just the version number fragment.\n  const fallback = body.search(\"2.\", { matchCase: true, matchWholeWord: false });\n  fallback.load(\"items\");\n  await context.sync();\n  if (fallback.items.length > 0) {\n    fallback.items[0].insertText(\"1.\", \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# \"Wireframes version 2.\" -> \"Wireframes version 1.\" style edit:\n# the paragraph's \"Version 2.\" becomes \"Version 1.\", with the two\n# existing runs (\"Versi\"/\"on\" and \" 2\"/\".\") each collapsing into a\n# single run (\"Version\" and \" 1.\"), while the \"_GoBack\" bookmark that\n# sits between \" 2\" and \".\" stays put, ending up right after \" 1.\".\n\n$d = $word.ActiveDocument\n\n# Locate the whole phrase first so the edit is anchored to the right\n# text wherever it lives in the document.\n$whole = $d.Content\n$whole.Find.ClearFormatting()\n$found = $whole.Find.Execute(\"Version 2.\")\nif (-not $found) {\n    throw \"Target text 'Version 2.' not found\"\n}\n$base = $whole.Start\n\n# --- \"Versi\" + \"on\" -> \"Version\" -------------------------------------\n# Office/Word only collapses adjacent runs into one when the assigned\n# text actually differs from the current contents, so round-trip\n# through a throwaway value first to force the merge even though the\n# final text (\"Version\") matches the original characters.\n$wordRange = $d.Range($base, $base + 7)\n$wordRange.Text = \"VersionTMP\"\n$wordRange2 = $d.Range($base, $base + 10)\n$wordRange2.Text = \"Version\"\n\n# --- \" 2\" -> \" 1\" (leave the bookmark and \".\" alone for now) ---------\n$numRange = $d.Range($base + 7, $base + 9)\n$numRange.Text = \" 1\"\n\n# --- remove the trailing \".\" run --------------------------------------\n$dotRange = $d.Range($base + 9, $base + 10)\n$dotRange.Delete()\n\n# --- re-attach \".\" right after \" 1\", ahead of the bookmark -----------\n$insertPoint = $d.Range($base + 9, $base + 9)\n$insertPoint.InsertBefore(\".\")\n"}
